$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.029.51'
$ws.Range("E2").Value = '  -1.80%  '
$ws.Range("D3").Value = '3.729.89'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''620.85'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '''180.19'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("D7").Value = '3.723.16'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.533'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("E11").Value = '  -4.53%  '
$ws.Range("D12").Value = '''0.485'
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").Value = '''40.58'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").Value = '''0.0000258'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '4.356.10'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '3.732.12'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '70.050.17'
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").Value = '''16.69'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").Value = '''505.39'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("E22").Value = '  -1.27%  '
$ws.Range("D23").Value = '''0.720'
$ws.Range("E23").Value = '  -4.54%  '
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").Value = '''86.61'
$ws.Range("E25").Value = '  -2.62%  '
$ws.Range("D26").Value = '''11.47'
$ws.Range("E26").Value = '  +2.85%  '
$ws.Range("D27").Value = '''13.09'
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("D28").Value = '''0.0000136'
$ws.Range("E28").Value = '  +20.60%  '
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").Value = '''2.48'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("D32").Value = '''7.92'
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("D33").Value = '''31.18'
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = '''6.16'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("D39").Value = '''0.339'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("E40").Value = '  -7.39%  '
$ws.Range("D41").Value = '''50.19'
$ws.Range("E41").Value = '  -3.11%  '
$ws.Range("D42").Value = '''45.51'
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("D43").Value = '''432.88'
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").Value = '''2.89'
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("D45").Value = '''8.68'
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = '2.997.60'
$ws.Range("E46").Value = '  -5.76%  '
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("D48").Value = '''27.47'
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("D50").Value = '''137.06'
$ws.Range("E50").Value = '  -2.48%  '
$ws.Range("E51").Value = '  +1.12%  '
